$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously had a 4x4 Sending-cluster x Target-cluster matrix
# (ECs/FAPs/MuSCs/Resolving-Mac) x (ECs/FAPs/MuSCs/Resolving-Mac).
# The updated TPM pipeline drops "Resolving-Mac" as a Target cluster, so the
# 4 rows that paired it as a target (old rows 14-17) are removed, leaving a
# 4x3 matrix (12 data rows) with refreshed NATMI scores throughout.
$ws.Rows.Item(14).EntireRow.Delete()
$ws.Rows.Item(14).EntireRow.Delete()
$ws.Rows.Item(14).EntireRow.Delete()
$ws.Rows.Item(14).EntireRow.Delete()

# Refresh data rows 2-13 (Sending cluster / Target cluster / scores) with the
# values recomputed from the new TPM input

# row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.893344666666666
$ws.Range("H2").Value = 5.680033999999999
$ws.Range("I2").Value = 0.05525983881677096
$ws.Range("J2").Value = 0.05525983881677096
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04214833333333334
$ws.Range("N2").Value = 0.126445
$ws.Range("O2").Value = 0.03198040784283177
$ws.Range("P2").Value = 0.03198040784283177
$ws.Range("Q2").Value = 0.07980132212555555
$ws.Range("R2").Value = 0.7182118991299999
$ws.Range("S2").Value = 0.001767232182689481
$ws.Range("T2").Value = 0.001767232182689481

# row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.893344666666666
$ws.Range("H3").Value = 5.680033999999999
$ws.Range("I3").Value = 0.05525983881677096
$ws.Range("J3").Value = 0.05525983881677096
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.30541
$ws.Range("N3").Value = 0.91623
$ws.Range("O3").Value = 0.2317324455521195
$ws.Range("P3").Value = 0.2317324455521195
$ws.Range("Q3").Value = 0.5782463946466666
$ws.Range("R3").Value = 5.204217551819999
$ws.Range("S3").Value = 0.01280549758982627
$ws.Range("T3").Value = 0.01280549758982627

# row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.893344666666666
$ws.Range("H4").Value = 5.680033999999999
$ws.Range("I4").Value = 0.05525983881677096
$ws.Range("J4").Value = 0.05525983881677096
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9703840000000001
$ws.Range("N4").Value = 2.911152
$ws.Range("O4").Value = 0.7362871466050487
$ws.Range("P4").Value = 0.7362871466050488
$ws.Range("Q4").Value = 1.837271371018667
$ws.Range("R4").Value = 16.535442339168
$ws.Range("S4").Value = 0.0406871090442552
$ws.Range("T4").Value = 0.0406871090442552

# row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.159773333333334
$ws.Range("H5").Value = 12.47932
$ws.Range("I5").Value = 0.1214086415227279
$ws.Range("J5").Value = 0.1214086415227279
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04214833333333334
$ws.Range("N5").Value = 0.126445
$ws.Range("O5").Value = 0.03198040784283177
$ws.Range("P5").Value = 0.03198040784283177
$ws.Range("Q5").Value = 0.1753275130444445
$ws.Range("R5").Value = 1.5779476174
$ws.Range("S5").Value = 0.003882697871540998
$ws.Range("T5").Value = 0.003882697871540998

# row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.159773333333334
$ws.Range("H6").Value = 12.47932
$ws.Range("I6").Value = 0.1214086415227279
$ws.Range("J6").Value = 0.1214086415227279
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.30541
$ws.Range("N6").Value = 0.91623
$ws.Range("O6").Value = 0.2317324455521195
$ws.Range("P6").Value = 0.2317324455521195
$ws.Range("Q6").Value = 1.270436373733334
$ws.Range("R6").Value = 11.4339273636
$ws.Range("S6").Value = 0.02813432141122233
$ws.Range("T6").Value = 0.02813432141122233

# row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.159773333333334
$ws.Range("H7").Value = 12.47932
$ws.Range("I7").Value = 0.1214086415227279
$ws.Range("J7").Value = 0.1214086415227279
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9703840000000001
$ws.Range("N7").Value = 2.911152
$ws.Range("O7").Value = 0.7362871466050487
$ws.Range("P7").Value = 0.7362871466050488
$ws.Range("Q7").Value = 4.036577486293334
$ws.Range("R7").Value = 36.32919737664001
$ws.Range("S7").Value = 0.08939162223996455
$ws.Range("T7").Value = 0.08939162223996457

# row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.240212333333333
$ws.Range("H8").Value = 3.720637
$ws.Range("I8").Value = 0.03619728348733726
$ws.Range("J8").Value = 0.03619728348733727
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04214833333333334
$ws.Range("N8").Value = 0.126445
$ws.Range("O8").Value = 0.03198040784283177
$ws.Range("P8").Value = 0.03198040784283177
$ws.Range("Q8").Value = 0.05227288282944444
$ws.Range("R8").Value = 0.470455945465
$ws.Range("S8").Value = 0.001157603888727645
$ws.Range("T8").Value = 0.001157603888727646

# row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.240212333333333
$ws.Range("H9").Value = 3.720637
$ws.Range("I9").Value = 0.03619728348733726
$ws.Range("J9").Value = 0.03619728348733727
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.30541
$ws.Range("N9").Value = 0.91623
$ws.Range("O9").Value = 0.2317324455521195
$ws.Range("P9").Value = 0.2317324455521195
$ws.Range("Q9").Value = 0.3787732487233333
$ws.Range("R9").Value = 3.40895923851
$ws.Range("S9").Value = 0.008388085024864015
$ws.Range("T9").Value = 0.008388085024864017

# row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.240212333333333
$ws.Range("H10").Value = 3.720637
$ws.Range("I10").Value = 0.03619728348733726
$ws.Range("J10").Value = 0.03619728348733727
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.9703840000000001
$ws.Range("N10").Value = 2.911152
$ws.Range("O10").Value = 0.7362871466050487
$ws.Range("P10").Value = 0.7362871466050488
$ws.Range("Q10").Value = 1.203482204869333
$ws.Range("R10").Value = 10.831339843824
$ws.Range("S10").Value = 0.0266515945737456
$ws.Range("T10").Value = 0.02665159457374561

# row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 26.96925
$ws.Range("H11").Value = 80.90774999999999
$ws.Range("I11").Value = 0.7871342361731639
$ws.Range("J11").Value = 0.7871342361731638
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04214833333333334
$ws.Range("N11").Value = 0.126445
$ws.Range("O11").Value = 0.03198040784283177
$ws.Range("P11").Value = 0.03198040784283177
$ws.Range("Q11").Value = 1.13670893875
$ws.Range("R11").Value = 10.23038044875
$ws.Range("S11").Value = 0.02517287389987364
$ws.Range("T11").Value = 0.02517287389987364

# row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 26.96925
$ws.Range("H12").Value = 80.90774999999999
$ws.Range("I12").Value = 0.7871342361731639
$ws.Range("J12").Value = 0.7871342361731638
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.30541
$ws.Range("N12").Value = 0.91623
$ws.Range("O12").Value = 0.2317324455521195
$ws.Range("P12").Value = 0.2317324455521195
$ws.Range("Q12").Value = 8.236678642499999
$ws.Range("R12").Value = 74.13010778249999
$ws.Range("S12").Value = 0.1824045415262069
$ws.Range("T12").Value = 0.1824045415262068

# row 13: Resolving-Mac -> MuSCs
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 26.96925
$ws.Range("H13").Value = 80.90774999999999
$ws.Range("I13").Value = 0.7871342361731639
$ws.Range("J13").Value = 0.7871342361731638
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.9703840000000001
$ws.Range("N13").Value = 2.911152
$ws.Range("O13").Value = 0.7362871466050487
$ws.Range("P13").Value = 0.7362871466050488
$ws.Range("Q13").Value = 26.170528692
$ws.Range("R13").Value = 235.534758228
$ws.Range("S13").Value = 0.5795568207470834
$ws.Range("T13").Value = 0.5795568207470834
